$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 174, shifting existing rows 174:264 down to 175:265
$ws.Rows.Item(174).Insert()

# Populate the newly inserted row 174 with the new record
$ws.Cells.Item(174, 1).Value = 8
$ws.Cells.Item(174, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(174, 3).Value = "Coquimbo"
$ws.Cells.Item(174, 4).Value = 44572
$ws.Cells.Item(174, 5).Value = 4
$ws.Cells.Item(174, 6).Value = 100114013
$ws.Cells.Item(174, 7).Value = "Zanahoria"
$ws.Cells.Item(174, 8).Value = "Sin especificar"
$ws.Cells.Item(174, 9).Value = "Primera"
$ws.Cells.Item(174, 10).Value = 600
$ws.Cells.Item(174, 11).Value = 5500
$ws.Cells.Item(174, 12).Value = 6000
$ws.Cells.Item(174, 13).Value = 5750
$ws.Cells.Item(174, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(174, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(174, 16).Value = 288
$ws.Cells.Item(174, 17).Value = 20
$ws.Cells.Item(174, 18).Value = "Hortaliza"

# Ensure the date cell keeps the date-style numeric formatting used by the rest of column D
$ws.Cells.Item(174, 4).NumberFormat = $ws.Cells.Item(175, 4).NumberFormat
